$d = $word.ActiveDocument

# The paragraph that ends in "...projekt_tasks_modify.html" needs to be
# split right after that filename, and a brand-new line needs to be
# inserted there:
#
#   <tab><tab><tab>Bejegyzések<tab><tab><tab><tab>- Kész<tab><tab><tab>project_tasks_posts.html
#
# followed by the rest of the original paragraph (unchanged).
#
# Using Find/Execute with wdReplaceOne and the special ^p / ^t codes lets
# Word re-flow the paragraph break and tab stops for us in one shot.

$old = "projekt_tasks_modify.html"
$new = "projekt_tasks_modify.html^p^t^t^tBejegyzések^t^t^t^t- Kész^t^t^tproject_tasks_posts.html"

$range = $d.Content
$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
